# cedulas procesadas la noche de 1/10/2025
# Appends the newly processed "cedula" numbers (rows 683-874, column A)
# to the Hoja1 sheet, reproducing the same cell styling (date / green /
# red / plain) used for the existing blocks of data above them, and
# updates the sheet selection to the newly added block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set values for rows 683-874 (column A) ---
$values = @{
  683 = 45931
  684 = 30080356
  685 = 21236224
  686 = 21229026
  687 = 39727147
  688 = 17329021
  689 = 39728361
  690 = 40367231
  691 = 37671025
  692 = 30080371
  693 = 40188151
  694 = 40334977
  695 = 1120498000
  696 = 40403787
  697 = 1121891980
  698 = 1121864632
  699 = 40332711
  700 = 52655858
  701 = 1069900532
  702 = 1121847839
  703 = 40413127
  704 = 40413367
  705 = 39648983
  706 = 1007702703
  707 = 1121857800
  708 = 40218053
  709 = 1006442026
  710 = 1006774169
  711 = 1121906325
  712 = 1121892958
  713 = 1072395724
  714 = 1121842098
  715 = 1121881992
  716 = 1121855560
  717 = 1001092252
  718 = 1143252992
  719 = 35261784
  720 = 1121964736
  721 = 1006775869
  722 = 40394987
  723 = 1075685562
  724 = 40442006
  725 = 1001118801
  726 = 26163691
  727 = 40390505
  728 = 1121955504
  729 = 1000832681
  730 = 1121849388
  731 = 1121944032
  732 = 1122652964
  733 = 1121818890
  734 = 1121954646
  735 = 21242719
  736 = 1121844971
  737 = 1006775707
  738 = 40215468
  739 = 1075277227
  740 = 1118291987
  741 = 40443728
  742 = 1121873859
  743 = 1121821756
  744 = 42031731
  745 = 1122653540
  746 = 30083822
  747 = 86086391
  748 = 40394841
  749 = 1006828058
  750 = 1121960138
  751 = 40334349
  752 = 20851383
  753 = 1121836742
  754 = 1121871447
  755 = 1007228432
  756 = 40218016
  757 = 1121824999
  758 = 1121917747
  759 = 1121827017
  760 = 1057574560
  761 = 1121862649
  762 = 40217183
  763 = 40325663
  764 = 1006798198
  765 = 1121852594
  766 = 1121946818
  767 = 1125552959
  768 = 40331222
  769 = 1056782595
  770 = 1006772734
  771 = 1121893657
  772 = 1193092381
  773 = 1121849200
  774 = 1012331970
  775 = 1012342436
  776 = 1121853595
  777 = 35587282
  778 = 1014184746
  779 = 1023035549
  780 = 40326656
  781 = 1007413869
  782 = 1006828141
  783 = 1121841888
  784 = 52619648
  785 = 1122138078
  786 = 40186281
  787 = 40371926
  788 = 1121896390
  789 = 35261879
  790 = 40401474
  791 = 1121824236
  792 = 1123561577
  793 = 40325807
  794 = 40441092
  795 = 40188434
  796 = 40402375
  797 = 40444722
  798 = 1193557701
  799 = 38290869
  800 = 52205911
  801 = 21231994
  802 = 1121920906
  803 = 42548918
  804 = 40441978
  805 = 1123801589
  806 = 1121817317
  807 = 1127391633
  808 = 1121915258
  809 = 1121899304
  810 = 40386471
  811 = 1121925708
  812 = 1121846671
  813 = 40400188
  814 = 1121968595
  815 = 1006874896
  816 = 51883053
  817 = 1030623737
  818 = 40388244
  819 = 37397372
  820 = 52977026
  821 = 40392197
  822 = 1121823721
  823 = 52862750
  824 = 1007741689
  825 = 40189602
  826 = 1234790567
  827 = 1006335084
  828 = 1121849738
  829 = 1099210766
  830 = 1121918538
  831 = 1116439136
  832 = 1234789507
  833 = 1033748818
  834 = 1121908839
  835 = 1121946821
  836 = 1121915532
  837 = 1120378354
  838 = 1006835763
  839 = 1007816099
  840 = 1005566462
  841 = 1005294945
  842 = 1122655308
  843 = 1123561702
  844 = 1006775023
  845 = 1006797846
  846 = 1003583504
  847 = 1121146010
  848 = 1122653301
  849 = 1121879243
  850 = 1120365722
  851 = 1121948797
  852 = 1120559199
  853 = 1006690427
  854 = 1121720012
  855 = 1121961367
  856 = 40186992
  857 = 1010123626
  858 = 1234789441
  859 = 1123431666
  860 = 1006859731
  861 = 1120385112
  862 = 1121921356
  863 = 1006774144
  864 = 1123140904
  865 = 1010031875
  866 = 46376887
  867 = 1006820624
  868 = 1121833506
  869 = 1121888058
  870 = 1122919487
  871 = 1193221023
  872 = 1121872788
  873 = 1069735886
  874 = 1007437551
}

foreach ($r in $values.Keys) {
  $ws.Cells.Item($r, 1).Value = $values[$r]
}

# --- Apply styles by copying format from an existing template cell of the same style ---
# style map: 0 -> A666 (no fill/default), 1 -> A1 (date number format),
#            3 -> A656 (green fill), 5 -> A657 (red fill)
$styleTemplates = @{ 0 = "A666"; 1 = "A1"; 3 = "A656"; 5 = "A657" }

$styleRuns = @(
  @{ Start = 683; End = 683; Style = 1 }
  @{ Start = 684; End = 684; Style = 5 }
  @{ Start = 685; End = 685; Style = 3 }
  @{ Start = 686; End = 691; Style = 5 }
  @{ Start = 692; End = 692; Style = 3 }
  @{ Start = 693; End = 699; Style = 5 }
  @{ Start = 700; End = 701; Style = 3 }
  @{ Start = 702; End = 702; Style = 5 }
  @{ Start = 703; End = 703; Style = 3 }
  @{ Start = 704; End = 706; Style = 5 }
  @{ Start = 707; End = 707; Style = 3 }
  @{ Start = 708; End = 709; Style = 5 }
  @{ Start = 710; End = 711; Style = 3 }
  @{ Start = 712; End = 713; Style = 5 }
  @{ Start = 714; End = 717; Style = 3 }
  @{ Start = 718; End = 719; Style = 5 }
  @{ Start = 720; End = 721; Style = 3 }
  @{ Start = 722; End = 722; Style = 5 }
  @{ Start = 723; End = 723; Style = 3 }
  @{ Start = 724; End = 725; Style = 5 }
  @{ Start = 726; End = 726; Style = 3 }
  @{ Start = 727; End = 730; Style = 5 }
  @{ Start = 731; End = 732; Style = 3 }
  @{ Start = 733; End = 736; Style = 5 }
  @{ Start = 737; End = 737; Style = 3 }
  @{ Start = 738; End = 739; Style = 5 }
  @{ Start = 740; End = 740; Style = 0 }
  @{ Start = 741; End = 741; Style = 3 }
  @{ Start = 742; End = 743; Style = 5 }
  @{ Start = 744; End = 744; Style = 3 }
  @{ Start = 745; End = 745; Style = 5 }
  @{ Start = 746; End = 746; Style = 3 }
  @{ Start = 747; End = 747; Style = 5 }
  @{ Start = 748; End = 748; Style = 3 }
  @{ Start = 749; End = 749; Style = 5 }
  @{ Start = 750; End = 753; Style = 3 }
  @{ Start = 754; End = 754; Style = 5 }
  @{ Start = 755; End = 755; Style = 3 }
  @{ Start = 756; End = 756; Style = 5 }
  @{ Start = 757; End = 757; Style = 3 }
  @{ Start = 758; End = 759; Style = 5 }
  @{ Start = 760; End = 760; Style = 3 }
  @{ Start = 761; End = 761; Style = 5 }
  @{ Start = 762; End = 763; Style = 3 }
  @{ Start = 764; End = 771; Style = 5 }
  @{ Start = 772; End = 772; Style = 3 }
  @{ Start = 773; End = 773; Style = 5 }
  @{ Start = 774; End = 777; Style = 3 }
  @{ Start = 778; End = 779; Style = 5 }
  @{ Start = 780; End = 780; Style = 3 }
  @{ Start = 781; End = 782; Style = 5 }
  @{ Start = 783; End = 784; Style = 3 }
  @{ Start = 785; End = 785; Style = 5 }
  @{ Start = 786; End = 787; Style = 3 }
  @{ Start = 788; End = 788; Style = 5 }
  @{ Start = 789; End = 789; Style = 3 }
  @{ Start = 790; End = 794; Style = 5 }
  @{ Start = 795; End = 798; Style = 3 }
  @{ Start = 799; End = 799; Style = 5 }
  @{ Start = 800; End = 804; Style = 3 }
  @{ Start = 805; End = 805; Style = 5 }
  @{ Start = 806; End = 806; Style = 3 }
  @{ Start = 807; End = 807; Style = 5 }
  @{ Start = 808; End = 809; Style = 3 }
  @{ Start = 810; End = 810; Style = 5 }
  @{ Start = 811; End = 816; Style = 3 }
  @{ Start = 817; End = 818; Style = 5 }
  @{ Start = 819; End = 819; Style = 3 }
  @{ Start = 820; End = 820; Style = 5 }
  @{ Start = 821; End = 822; Style = 3 }
  @{ Start = 823; End = 827; Style = 5 }
  @{ Start = 828; End = 831; Style = 3 }
  @{ Start = 832; End = 832; Style = 5 }
  @{ Start = 833; End = 833; Style = 0 }
  @{ Start = 834; End = 834; Style = 3 }
  @{ Start = 835; End = 837; Style = 5 }
  @{ Start = 838; End = 838; Style = 3 }
  @{ Start = 839; End = 842; Style = 5 }
  @{ Start = 843; End = 843; Style = 3 }
  @{ Start = 844; End = 848; Style = 5 }
  @{ Start = 849; End = 855; Style = 3 }
  @{ Start = 856; End = 874; Style = 5 }
)

foreach ($run in $styleRuns) {
  $template = $styleTemplates[$run.Style]
  $ws.Range($template).Copy()
  $destRange = $ws.Range("A$($run.Start):A$($run.End)")
  $destRange.PasteSpecial(-4122)
}

# --- Update sheet view: selection + scroll position ---
$ws.Range("A684:A874").Select()

# Best-effort: move the visible top-left cell down to the new block
# (harmless no-op on hosts that don't persist window scroll state).
try {
  $excel.ActiveWindow.ScrollRow = 669
  $excel.ActiveWindow.ScrollColumn = 1
} catch {}
